# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.553.87'
$ws.Range('E2').Value = '  -1.27%  '
$ws.Range('D3').Value = '2.585.05'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.73'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.61'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.23%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -1.52%  '
$ws.Range('D9').Value = '2.585.67'
$ws.Range('E9').Value = '  -1.88%  '
$ws.Range('E10').Value = '  -4.03%  '
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('E12').Value = '  -1.77%  '
$ws.Range('E13').Value = '  -1.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.74'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.12%  '
$ws.Range('D15').Value = '3.056.35'
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('E16').Value = '  -2.78%  '
$ws.Range('D17').Value = '66.305.84'
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('D18').Value = '2.580.31'
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('E19').Value = '  -6.23%  '
$ws.Range('E20').Value = '  -4.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '352.05'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.22'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.61'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.25%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  -4.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '68.79'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.96'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -8.94%  '
$ws.Range('D28').Value = '2.716.98'
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('D29').Value = '0.0₃0987'
$ws.Range('E29').Value = '  -3.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '533.01'
$ws.Range('D30').ClearFormats()
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('E32').Value = '  -3.35%  '
$ws.Range('E33').Value = '  -3.15%  '
$ws.Range('E34').Value = '  -3.04%  '
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('E36').Value = '  -3.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '156.94'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('E38').Value = '  -2.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.360'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.01%  '
$ws.Range('E40').Value = '  +1.77%  '
$ws.Range('E41').Value = '  -1.58%  '
$ws.Range('E42').Value = '  -2.23%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  -2.79%  '
$ws.Range('D45').Value = '0.0₆0287'
$ws.Range('E45').Value = '  -4.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '149.28'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.566'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.93%  '
$ws.Range('E48').Value = '  -2.68%  '
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.597'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.64%  '
